$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns are treated as plain text so
# values like "7.30", "1.00", "0.0000104" keep their exact string form instead
# of being auto-converted to numbers (which would drop formatting/precision).
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '42.656.55'
$ws.Range("E2").Value = '  -0.70%  '

# Row 3
$ws.Range("D3").Value = '2.294.55'
$ws.Range("E3").Value = '  -0.59%  '

# Row 4
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.32%  '

# Row 5
$ws.Range("D5").Value = '311.65'
$ws.Range("E5").Value = '  -3.29%  '

# Row 6
$ws.Range("D6").Value = '103.39'
$ws.Range("E6").Value = '  -2.01%  '

# Row 7
$ws.Range("D7").Value = '0.624'
$ws.Range("E7").Value = '  -1.27%  '

# Row 8
$ws.Range("E8").Value = '  -0.02%  '

# Row 9
$ws.Range("D9").Value = '0.604'
$ws.Range("E9").Value = '  -1.26%  '

# Row 10
$ws.Range("D10").Value = '39.14'
$ws.Range("E10").Value = '  -3.62%  '

# Row 11
$ws.Range("E11").Value = '  -1.28%  '

# Row 12
$ws.Range("D12").Value = '8.25'
$ws.Range("E12").Value = '  -4.37%  '

# Row 13
$ws.Range("E13").Value = '  +0.77%  '

# Row 14
$ws.Range("D14").Value = '0.985'
$ws.Range("E14").Value = '  +0.73%  '

# Row 15
$ws.Range("D15").Value = '15.38'
$ws.Range("E15").Value = '  -0.12%  '

# Row 16
$ws.Range("D16").Value = '2.638.46'
$ws.Range("E16").Value = '  -0.72%  '

# Row 17
$ws.Range("D17").Value = '2.295.00'
$ws.Range("E17").Value = '  -0.16%  '

# Row 18
$ws.Range("D18").Value = '42.596.23'
$ws.Range("E18").Value = '  -0.35%  '

# Row 19
$ws.Range("D19").Value = '7.30'
$ws.Range("E19").Value = '  -3.54%  '

# Row 20
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.0000104'
$ws.Range("E20").Value = '  -1.68%  '

# Row 21
$ws.Range("B21").Value = 'InternetComputer(DFINITY)'
$ws.Range("C21").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D21").Value = '13.53'
$ws.Range("E21").Value = '  +0.43%  '

# Row 22
$ws.Range("D22").Value = '73.43'
$ws.Range("E22").Value = '  -0.64%  '

# Row 23
$ws.Range("D23").Value = '267.93'
$ws.Range("E23").Value = '  -1.68%  '

# Row 24
$ws.Range("D24").Value = '3.42'
$ws.Range("E24").Value = '  -4.92%  '

# Row 25
$ws.Range("E25").Value = '  -2.54%  '

# Row 26
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.30%  '

# Row 27
$ws.Range("D27").Value = '10.80'
$ws.Range("E27").Value = '  -1.58%  '

# Row 28
$ws.Range("D28").Value = '7.22'
$ws.Range("E28").Value = '  +16.58%  '

# Row 29
$ws.Range("E29").Value = '  -1.40%  '

# Row 30
$ws.Range("D30").Value = '22.38'
$ws.Range("E30").Value = '  -1.63%  '

# Row 31
$ws.Range("D31").Value = '35.93'
$ws.Range("E31").Value = '  -6.35%  '

# Row 32
$ws.Range("D32").Value = '164.65'
$ws.Range("E32").Value = '  -0.75%  '

# Row 33
$ws.Range("D33").Value = '0.0853'
$ws.Range("E33").Value = '  -4.06%  '

# Row 34
$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").Value = '2.63'
$ws.Range("E34").Value = '  +3.68%  '

# Row 35
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").Value = '0.130'
$ws.Range("E35").Value = '  -2.21%  '

# Row 36
$ws.Range("D36").Value = '0.112'
$ws.Range("E36").Value = '  -3.64%  '

# Row 37
$ws.Range("D37").Value = '4.54'
$ws.Range("E37").Value = '  -2.36%  '

# Row 38
$ws.Range("D38").Value = '0.0347'
$ws.Range("E38").Value = '  -2.79%  '

# Row 39
$ws.Range("D39").Value = '2.81'
$ws.Range("E39").Value = '  +2.31%  '

# Row 40
$ws.Range("D40").Value = '3.63'
$ws.Range("E40").Value = '  -3.10%  '

# Row 41
$ws.Range("D41").Value = '107.11'
$ws.Range("E41").Value = '  +5.09%  '

# Row 42
$ws.Range("D42").Value = '1.58'
$ws.Range("E42").Value = '  +0.45%  '

# Row 43
$ws.Range("D43").Value = '70.59'
$ws.Range("E43").Value = '  -0.44%  '

# Row 44
$ws.Range("E44").Value = '  +0.25%  '

# Row 45
$ws.Range("E45").Value = '  +0.02%  '

# Row 46
$ws.Range("D46").Value = '1.746.76'
$ws.Range("E46").Value = '  +9.43%  '

# Row 47
$ws.Range("D47").Value = '12.09'
$ws.Range("E47").Value = '  -3.13%  '

# Row 48
$ws.Range("D48").Value = '110.40'
$ws.Range("E48").Value = '  -3.46%  '

# Row 49
$ws.Range("D49").Value = '77.48'
$ws.Range("E49").Value = '  -6.62%  '

# Row 50
$ws.Range("E50").Value = '  -3.27%  '

# Row 51
$ws.Range("D51").Value = '8.64'
$ws.Range("E51").Value = '  -3.20%  '
